$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 397.16666  # H33 782.3333 -> 397.16666
$ws.Cells.Item(33, 9).Value = 397.16666  # I33 782.3333 -> 397.16666
$ws.Cells.Item(33, 11).Value = 397.16666  # K33 782.3333 -> 397.16666
$ws.Cells.Item(33, 13).Value = -168.16666  # M33 -553.3333 -> -168.16666
$ws.Cells.Item(55, 8).Value = 757.05884  # H55 726.1111 -> 757.05884
$ws.Cells.Item(55, 10).Value = 846.9091  # J55 793 -> 846.9091
$ws.Cells.Item(55, 12).Value = 846.9091  # L55 793 -> 846.9091
$ws.Cells.Item(55, 14).Value = -1274.9091  # N55 -1221 -> -1274.9091
$ws.Cells.Item(62, 8).Value = 9207.5  # H62 9991.5 -> 9207.5
$ws.Cells.Item(62, 9).Value = 3827.5  # I62 3950 -> 3827.5
$ws.Cells.Item(62, 10).Value = 11000.833  # J62 11199.8 -> 11000.833
$ws.Cells.Item(62, 11).Value = 3827.5  # K62 3950 -> 3827.5
$ws.Cells.Item(62, 12).Value = 11000.833  # L62 11199.8 -> 11000.833
$ws.Cells.Item(62, 13).Value = -3203.5  # M62 -3326 -> -3203.5
$ws.Cells.Item(62, 14).Value = -12248.833  # N62 -12447.8 -> -12248.833
$ws.Cells.Item(64, 8).Value = 0  # H64 1500 -> 0
$ws.Cells.Item(64, 9).Value = 0  # I64 1500 -> 0
$ws.Cells.Item(64, 11).Value = 0  # K64 1500 -> 0
$ws.Cells.Item(64, 13).ClearContents()  # M64 was -1252
$ws.Cells.Item(65, 8).Value = 9207.5  # H65 9991.5 -> 9207.5
$ws.Cells.Item(65, 9).Value = 3827.5  # I65 3950 -> 3827.5
$ws.Cells.Item(65, 10).Value = 11000.833  # J65 11199.8 -> 11000.833
$ws.Cells.Item(65, 11).Value = 19137.5  # K65 19750 -> 19137.5
$ws.Cells.Item(65, 12).Value = 55004.165  # L65 55999 -> 55004.165
$ws.Cells.Item(65, 13).Value = -16017.5  # M65 -16630 -> -16017.5
$ws.Cells.Item(65, 14).Value = -61244.165  # N65 -62239 -> -61244.165
$ws.Cells.Item(67, 8).Value = 0  # H67 1500 -> 0
$ws.Cells.Item(67, 9).Value = 0  # I67 1500 -> 0
$ws.Cells.Item(67, 11).Value = 0  # K67 1500 -> 0
$ws.Cells.Item(67, 13).ClearContents()  # M67 was -642
$ws.Cells.Item(98, 8).Value = 659.8  # H98 550 -> 659.8
$ws.Cells.Item(98, 9).Value = 324.75  # I98 550 -> 324.75
$ws.Cells.Item(98, 10).Value = 2000  # J98 0 -> 2000
$ws.Cells.Item(98, 11).Value = 324.75  # K98 550 -> 324.75
$ws.Cells.Item(98, 12).Value = 2000  # L98 0 -> 2000
$ws.Cells.Item(98, 13).Value = 1173.25  # M98 948 -> 1173.25
$ws.Cells.Item(98, 14).Value = -4996  # N98 None -> -4996
$ws.Cells.Item(107, 8).Value = 858.36365  # H107 756.6923 -> 858.36365
$ws.Cells.Item(107, 9).Value = 1034.5714  # I107 848.55554 -> 1034.5714
$ws.Cells.Item(107, 11).Value = 1034.5714  # K107 848.55554 -> 1034.5714
$ws.Cells.Item(107, 13).Value = 885.4286  # M107 1071.44446 -> 885.4286
$ws.Cells.Item(122, 8).Value = 659.8  # H122 550 -> 659.8
$ws.Cells.Item(122, 9).Value = 324.75  # I122 550 -> 324.75
$ws.Cells.Item(122, 10).Value = 2000  # J122 0 -> 2000
$ws.Cells.Item(122, 11).Value = 974.25  # K122 1650 -> 974.25
$ws.Cells.Item(122, 12).Value = 6000  # L122 0 -> 6000
$ws.Cells.Item(122, 13).Value = 1475.75  # M122 800 -> 1475.75
$ws.Cells.Item(122, 14).Value = -10900  # N122 None -> -10900
$ws.Cells.Item(137, 8).Value = 1981.5454  # H137 2154.4 -> 1981.5454
$ws.Cells.Item(137, 9).Value = 649.5  # I137 728.8 -> 649.5
$ws.Cells.Item(137, 11).Value = 1948.5  # K137 2186.4 -> 1948.5
$ws.Cells.Item(137, 13).Value = 601.5  # M137 363.6000000000004 -> 601.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(38, 8).Value = 13589.833  # H38 18254.5 -> 13589.833
$ws.Cells.Item(38, 9).Value = 3173  # I38 998 -> 3173
$ws.Cells.Item(38, 11).Value = 3173  # K38 998 -> 3173
$ws.Cells.Item(38, 13).Value = -2706  # M38 -531 -> -2706
$ws.Cells.Item(45, 8).Value = 3721.875  # H45 3631.111 -> 3721.875
$ws.Cells.Item(45, 9).Value = 2861  # I45 2788.6 -> 2861
$ws.Cells.Item(45, 10).Value = 4582.75  # J45 4684.25 -> 4582.75
$ws.Cells.Item(45, 11).Value = 2861  # K45 2788.6 -> 2861
$ws.Cells.Item(45, 12).Value = 4582.75  # L45 4684.25 -> 4582.75
$ws.Cells.Item(45, 13).Value = -2484  # M45 -2411.6 -> -2484
$ws.Cells.Item(45, 14).Value = -5336.75  # N45 -5438.25 -> -5336.75
$ws.Cells.Item(61, 8).Value = 4433.2  # H61 4504.8887 -> 4433.2
$ws.Cells.Item(61, 9).Value = 4333.2856  # I61 4424.1665 -> 4333.2856
$ws.Cells.Item(61, 11).Value = 4333.2856  # K61 4424.1665 -> 4333.2856
$ws.Cells.Item(61, 13).Value = -4121.2856  # M61 -4212.1665 -> -4121.2856
$ws.Cells.Item(74, 8).Value = 9340.375  # H74 14472.4 -> 9340.375
$ws.Cells.Item(74, 9).Value = 10560.429  # I74 14472.4 -> 10560.429
$ws.Cells.Item(74, 10).Value = 800  # J74 0 -> 800
$ws.Cells.Item(74, 11).Value = 10560.429  # K74 14472.4 -> 10560.429
$ws.Cells.Item(74, 12).Value = 800  # L74 0 -> 800
$ws.Cells.Item(74, 13).Value = -9686.429  # M74 -13598.4 -> -9686.429
$ws.Cells.Item(74, 14).Value = -2548  # N74 None -> -2548
$ws.Cells.Item(77, 8).Value = 9340.375  # H77 14472.4 -> 9340.375
$ws.Cells.Item(77, 9).Value = 10560.429  # I77 14472.4 -> 10560.429
$ws.Cells.Item(77, 10).Value = 800  # J77 0 -> 800
$ws.Cells.Item(77, 11).Value = 52802.145  # K77 72362 -> 52802.145
$ws.Cells.Item(77, 12).Value = 4000  # L77 0 -> 4000
$ws.Cells.Item(77, 13).Value = -48434.145  # M77 -67994 -> -48434.145
$ws.Cells.Item(77, 14).Value = -12736  # N77 None -> -12736
$ws.Cells.Item(102, 8).Value = 3869.158  # H102 3885.2632 -> 3869.158
$ws.Cells.Item(102, 9).Value = 1884.9231  # I102 1908.4615 -> 1884.9231
$ws.Cells.Item(102, 11).Value = 1884.9231  # K102 1908.4615 -> 1884.9231
$ws.Cells.Item(102, 13).Value = -262.9231  # M102 -286.4614999999999 -> -262.9231
$ws.Cells.Item(132, 8).Value = 4622  # H132 3716.8462 -> 4622
$ws.Cells.Item(132, 9).Value = 5521.375  # I132 4547 -> 5521.375
$ws.Cells.Item(132, 10).Value = 1024.5  # J132 949.6667 -> 1024.5
$ws.Cells.Item(132, 11).Value = 16564.125  # K132 13641 -> 16564.125
$ws.Cells.Item(132, 12).Value = 3073.5  # L132 2849.0001 -> 3073.5
$ws.Cells.Item(132, 13).Value = -14034.125  # M132 -11111 -> -14034.125
$ws.Cells.Item(132, 14).Value = -8133.5  # N132 -7909.0001 -> -8133.5
$ws.Cells.Item(136, 8).Value = 4433.2  # H136 4504.8887 -> 4433.2
$ws.Cells.Item(136, 9).Value = 4333.2856  # I136 4424.1665 -> 4333.2856
$ws.Cells.Item(136, 11).Value = 12999.8568  # K136 13272.4995 -> 12999.8568
$ws.Cells.Item(136, 13).Value = -10449.8568  # M136 -10722.4995 -> -10449.8568
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 555.8  # H80 557.73334 -> 555.8
$ws.Cells.Item(80, 10).Value = 1119.4  # J80 1125.2 -> 1119.4
$ws.Cells.Item(80, 12).Value = 1119.4  # L80 1125.2 -> 1119.4
$ws.Cells.Item(80, 14).Value = -3115.4  # N80 -3121.2 -> -3115.4
$ws.Cells.Item(83, 8).Value = 555.8  # H83 557.73334 -> 555.8
$ws.Cells.Item(83, 10).Value = 1119.4  # J83 1125.2 -> 1119.4
$ws.Cells.Item(83, 12).Value = 5597  # L83 5626 -> 5597
$ws.Cells.Item(83, 14).Value = -15581  # N83 -15610 -> -15581
$ws.Cells.Item(86, 8).Value = 5285.4  # H86 5899.154 -> 5285.4
$ws.Cells.Item(86, 9).Value = 3498.4285  # I86 3798.8 -> 3498.4285
$ws.Cells.Item(86, 10).Value = 6849  # J86 7211.875 -> 6849
$ws.Cells.Item(86, 11).Value = 3498.4285  # K86 3798.8 -> 3498.4285
$ws.Cells.Item(86, 12).Value = 6849  # L86 7211.875 -> 6849
$ws.Cells.Item(86, 13).Value = -2375.4285  # M86 -2675.8 -> -2375.4285
$ws.Cells.Item(86, 14).Value = -9095  # N86 -9457.875 -> -9095
$ws.Cells.Item(89, 8).Value = 5285.4  # H89 5899.154 -> 5285.4
$ws.Cells.Item(89, 9).Value = 3498.4285  # I89 3798.8 -> 3498.4285
$ws.Cells.Item(89, 10).Value = 6849  # J89 7211.875 -> 6849
$ws.Cells.Item(89, 11).Value = 17492.1425  # K89 18994 -> 17492.1425
$ws.Cells.Item(89, 12).Value = 34245  # L89 36059.375 -> 34245
$ws.Cells.Item(89, 13).Value = -11876.1425  # M89 -13378 -> -11876.1425
$ws.Cells.Item(89, 14).Value = -45477  # N89 -47291.375 -> -45477
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 529.8  # H105 403.5 -> 529.8
$ws.Cells.Item(105, 9).Value = 84.666664  # I105 138 -> 84.666664
$ws.Cells.Item(105, 10).Value = 1197.5  # J105 1200 -> 1197.5
$ws.Cells.Item(105, 11).Value = 84.666664  # K105 138 -> 84.666664
$ws.Cells.Item(105, 12).Value = 1197.5  # L105 1200 -> 1197.5
$ws.Cells.Item(105, 13).Value = 1662.333336  # M105 1609 -> 1662.333336
$ws.Cells.Item(105, 14).Value = -4691.5  # N105 -4694 -> -4691.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1242.2222  # H5 1297.7 -> 1242.2222
$ws.Cells.Item(5, 9).Value = 965.1667  # I5 981.6667 -> 965.1667
$ws.Cells.Item(5, 10).Value = 1796.3334  # J5 1771.75 -> 1796.3334
$ws.Cells.Item(5, 11).Value = 2895.5001  # K5 2945.0001 -> 2895.5001
$ws.Cells.Item(5, 12).Value = 5389.0002  # L5 5315.25 -> 5389.0002
$ws.Cells.Item(5, 13).Value = -2783.5001  # M5 -2833.0001 -> -2783.5001
$ws.Cells.Item(5, 14).Value = -5613.0002  # N5 -5539.25 -> -5613.0002
$ws.Cells.Item(135, 8).Value = 1242.2222  # H135 1297.7 -> 1242.2222
$ws.Cells.Item(135, 9).Value = 965.1667  # I135 981.6667 -> 965.1667
$ws.Cells.Item(135, 10).Value = 1796.3334  # J135 1771.75 -> 1796.3334
$ws.Cells.Item(135, 11).Value = 8686.5003  # K135 8835.0003 -> 8686.5003
$ws.Cells.Item(135, 12).Value = 16167.0006  # L135 15945.75 -> 16167.0006
$ws.Cells.Item(135, 13).Value = -6151.5003  # M135 -6300.0003 -> -6151.5003
$ws.Cells.Item(135, 14).Value = -21237.0006  # N135 -21015.75 -> -21237.0006
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 20014800  # H7 31111112 -> 20014800
$ws.Cells.Item(7, 9).Value = 20022222  # I7 22857142 -> 20022222
$ws.Cells.Item(7, 10).Value = 20003666  # J7 60000000 -> 20003666
$ws.Cells.Item(7, 11).Value = 20022222  # K7 22857142 -> 20022222
$ws.Cells.Item(7, 12).Value = 20003666  # L7 60000000 -> 20003666
$ws.Cells.Item(7, 13).Value = -20022110  # M7 -22857030 -> -20022110
$ws.Cells.Item(7, 14).Value = -20003890  # N7 -60000224 -> -20003890
$ws.Cells.Item(8, 8).Value = 20014800  # H8 31111112 -> 20014800
$ws.Cells.Item(8, 9).Value = 20022222  # I8 22857142 -> 20022222
$ws.Cells.Item(8, 10).Value = 20003666  # J8 60000000 -> 20003666
$ws.Cells.Item(8, 11).Value = 20022222  # K8 22857142 -> 20022222
$ws.Cells.Item(8, 12).Value = 20003666  # L8 60000000 -> 20003666
$ws.Cells.Item(8, 13).Value = -20022083  # M8 -22857003 -> -20022083
$ws.Cells.Item(8, 14).Value = -20003944  # N8 -60000278 -> -20003944
$ws.Cells.Item(14, 8).Value = 57005  # H14 500 -> 57005
$ws.Cells.Item(14, 9).Value = 0  # I14 500 -> 0
$ws.Cells.Item(14, 10).Value = 57005  # J14 0 -> 57005
$ws.Cells.Item(14, 11).Value = 0  # K14 500 -> 0
$ws.Cells.Item(14, 12).Value = 57005  # L14 0 -> 57005
$ws.Cells.Item(14, 13).ClearContents()  # M14 was -332
$ws.Cells.Item(14, 14).Value = -57341  # N14 None -> -57341
$ws.Cells.Item(36, 8).Value = 1988  # H36 2858.1428 -> 1988
$ws.Cells.Item(36, 9).Value = 3508.5  # I36 1379.25 -> 3508.5
$ws.Cells.Item(36, 10).Value = 1481.1666  # J36 4830 -> 1481.1666
$ws.Cells.Item(36, 11).Value = 3508.5  # K36 1379.25 -> 3508.5
$ws.Cells.Item(36, 12).Value = 1481.1666  # L36 4830 -> 1481.1666
$ws.Cells.Item(36, 13).Value = -3023.5  # M36 -894.25 -> -3023.5
$ws.Cells.Item(36, 14).Value = -2451.1666  # N36 -5800 -> -2451.1666
$ws.Cells.Item(43, 8).Value = 19629.3  # H43 14209.214 -> 19629.3
$ws.Cells.Item(43, 9).Value = 0  # I43 734.75 -> 0
$ws.Cells.Item(43, 10).Value = 19629.3  # J43 19599 -> 19629.3
$ws.Cells.Item(43, 11).Value = 0  # K43 734.75 -> 0
$ws.Cells.Item(43, 12).Value = 19629.3  # L43 19599 -> 19629.3
$ws.Cells.Item(43, 13).ClearContents()  # M43 was -583.75
$ws.Cells.Item(43, 14).Value = -19931.3  # N43 -19901 -> -19931.3
$ws.Cells.Item(102, 8).Value = 3130.2  # H102 2927.4546 -> 3130.2
$ws.Cells.Item(102, 9).Value = 2398.8572  # I102 2211.5 -> 2398.8572
$ws.Cells.Item(102, 11).Value = 2398.8572  # K102 2211.5 -> 2398.8572
$ws.Cells.Item(102, 13).Value = -776.8571999999999  # M102 -589.5 -> -776.8571999999999
$ws.Cells.Item(126, 8).Value = 4772.857  # H126 5151.6665 -> 4772.857
$ws.Cells.Item(126, 9).Value = 3682  # I126 3977.5 -> 3682
$ws.Cells.Item(126, 11).Value = 11046  # K126 11932.5 -> 11046
$ws.Cells.Item(126, 13).Value = -8576  # M126 -9462.5 -> -8576
$ws.Cells.Item(135, 8).Value = 101972.5  # H135 94296.664 -> 101972.5
$ws.Cells.Item(135, 10).Value = 101972.5  # J135 94296.664 -> 101972.5
$ws.Cells.Item(135, 12).Value = 101972.5  # L135 94296.664 -> 101972.5
$ws.Cells.Item(135, 14).Value = -112112.5  # N135 -104436.664 -> -112112.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 6091  # H40 6645.5293 -> 6091
$ws.Cells.Item(40, 9).Value = 5732.9375  # I40 6355.143 -> 5732.9375
$ws.Cells.Item(40, 11).Value = 5732.9375  # K40 6355.143 -> 5732.9375
$ws.Cells.Item(40, 13).Value = -5596.9375  # M40 -6219.143 -> -5596.9375
$ws.Cells.Item(122, 8).Value = 5000  # H122 0 -> 5000
$ws.Cells.Item(122, 9).Value = 5000  # I122 0 -> 5000
$ws.Cells.Item(122, 11).Value = 15000  # K122 0 -> 15000
$ws.Cells.Item(122, 13).Value = -12550  # M122 None -> -12550
$ws.Cells.Item(132, 8).Value = 2585.6667  # H132 1749.75 -> 2585.6667
$ws.Cells.Item(132, 10).Value = 4257.5  # J132 0 -> 4257.5
$ws.Cells.Item(132, 12).Value = 12772.5  # L132 0 -> 12772.5
$ws.Cells.Item(132, 14).Value = -17832.5  # N132 None -> -17832.5
$ws.Cells.Item(136, 8).Value = 1347.5  # H136 1297 -> 1347.5
$ws.Cells.Item(136, 9).Value = 1347.5  # I136 1273.75 -> 1347.5
$ws.Cells.Item(136, 10).Value = 0  # J136 1390 -> 0
$ws.Cells.Item(136, 11).Value = 4042.5  # K136 3821.25 -> 4042.5
$ws.Cells.Item(136, 12).Value = 0  # L136 4170 -> 0
$ws.Cells.Item(136, 13).Value = -1492.5  # M136 -1271.25 -> -1492.5
$ws.Cells.Item(136, 14).ClearContents()  # N136 was -9270
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 51500  # H64 53000 -> 51500
$ws.Cells.Item(64, 10).Value = 51500  # J64 53000 -> 51500
$ws.Cells.Item(64, 12).Value = 51500  # L64 53000 -> 51500
$ws.Cells.Item(64, 14).Value = -51996  # N64 -53496 -> -51996
$ws.Cells.Item(67, 8).Value = 51500  # H67 53000 -> 51500
$ws.Cells.Item(67, 10).Value = 51500  # J67 53000 -> 51500
$ws.Cells.Item(67, 12).Value = 51500  # L67 53000 -> 51500
$ws.Cells.Item(67, 14).Value = -53216  # N67 -54716 -> -53216
$ws.Cells.Item(122, 8).Value = 2000  # H122 0 -> 2000
$ws.Cells.Item(122, 10).Value = 2000  # J122 0 -> 2000
$ws.Cells.Item(122, 12).Value = 6000  # L122 0 -> 6000
$ws.Cells.Item(122, 14).Value = -10900  # N122 None -> -10900
$ws.Cells.Item(132, 8).Value = 1921.7693  # H132 2025.7273 -> 1921.7693
$ws.Cells.Item(132, 9).Value = 1978.3  # I132 2064.7778 -> 1978.3
$ws.Cells.Item(132, 10).Value = 1733.3334  # J132 1850 -> 1733.3334
$ws.Cells.Item(132, 11).Value = 5934.9  # K132 6194.3334 -> 5934.9
$ws.Cells.Item(132, 12).Value = 5200.0002  # L132 5550 -> 5200.0002
$ws.Cells.Item(132, 13).Value = -3404.9  # M132 -3664.3334 -> -3404.9
$ws.Cells.Item(132, 14).Value = -10260.0002  # N132 -10610 -> -10260.0002
